{"js": "// Apply the \"week 1 planning\" text edits to Activity1.docx\n// Each entry is an exact-match search over the whole paragraph text followed\n// by a full replacement of that run's text (search uses matchCase to avoid\n// unwanted partial matches).\nconst replacements = [\n  {\n    find: \"We have been tasked by Waypoint to investigate inequity in Iowa potentially exhibited in the services they offer and evictions in the state. Let\\u2019s start by understanding inequity. Afterwards you will research on Iowa\\u2019s history with racial inequities and homelessness.\",\n    replace: \"Waypoint Services has tasked us to investigate a community service to help people who are houseless. Let\\u2019s start by understanding inequity in housing.\"\n  },\n  {\n    find: \"As part of the first activity, please create a shared Google Drive folder for your group this block and upload your .docx file to that folder. Once uploaded you will all be able to work on that document at one time.\",\n    replace: \"As part of the first activity, please create a shared Google Drive folder for your group this block and upload your .docx file to that folder. Once uploaded, you will all be able to work on that document at one time.\"\n  },\n  {\n    find: \"In this activity you will all be responsible for the searches and reviewing of items founds. Generally you will all have assigned roles but for in this case you should all be working on the shared document adding different items as explained below. I would like someone to additionally be chosen as a facilitator.\",\n    replace: \"In this activity, you will all be responsible for conducting searches and reviewing items found. Generally, you will all have assigned roles, but in this case, you should all be working on the shared document, adding different items as explained below. Please choose a facilitator.\"\n  },\n  {\n    find: \"for this question. In your group I would like you to discuss among yourselves how you would define inequity. Write a one sentence definition.\",\n    replace: \"for this question. In your group, I would like you to discuss among yourselves how you would define inequity. Write a one-sentence definition.\"\n  },\n  {\n    find: \"Now you can start to use the internet. Each group member should spend a few minutes looking for different definitions of inequity and examples. Each group member should select a different example. Considering looking for something you have a personal connection to or have thought about before. Include your examples below (make sure to cite your sources, websites are fine in this case).\",\n    replace: \"Now you can start to use the internet. Each group member should spend a few minutes looking for different definitions of inequity and examples. Each group member should select a different example. Consider looking for something you have a personal connection to or have thought about before. Include your examples below (make sure to cite your sources; websites are fine in this case).\"\n  },\n  {\n    find: \"Now that you have a definition for inequity you can start too look for it in Iowa. In groups I want you to look for some anecdotal evidence. This type of evidence is primarily what Waypoint has founded their suspicions of inequity on. We are specifically looking for inequity in Iowa and how it relates to having a home.\",\n    replace: \"Now that you have a definition for inequity, you can start to look for it in Iowa. In groups, please look for some anecdotal evidence. This type of evidence is primarily what Waypoint has founded their suspicions of inequity on. We are specifically looking for inequity in Iowa and how it relates to having a home.\"\n  }\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + find);\n  }\n\n  for (const item of results.items) {\n    item.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"week 1 planning\" text edits to Activity1.docx\n#\n# Note: using Find.Execute(..., Replace:=wdReplaceOne) directly collapses the\n# matched range's run together with any immediately-adjacent run that shares\n# identical formatting (e.g. the single-space run that separates a bold\n# label like \"Google Drive:\" from the sentence that follows). That would\n# change the run layout beyond what the source edit touches. Instead, locate\n# the target text with Find.Execute() (no replacement arguments), then\n# Delete() + InsertAfter() on the now-narrowed range, which replaces only\n# that range's text and leaves neighboring runs untouched.\n$d = $word.ActiveDocument\n$apos = [char]0x2019\n\nfunction Replace-Text($oldText, $newText) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"Replacement failed for: $oldText\"\n    }\n    $rng.Delete()\n    $rng.InsertAfter($newText)\n}\n\nReplace-Text `\n    \"We have been tasked by Waypoint to investigate inequity in Iowa potentially exhibited in the services they offer and evictions in the state. Let${apos}s start by understanding inequity. Afterwards you will research on Iowa${apos}s history with racial inequities and homelessness.\" `\n    \"Waypoint Services has tasked us to investigate a community service to help people who are houseless. Let${apos}s start by understanding inequity in housing.\"\n\nReplace-Text `\n    \"As part of the first activity, please create a shared Google Drive folder for your group this block and upload your .docx file to that folder. Once uploaded you will all be able to work on that document at one time.\" `\n    \"As part of the first activity, please create a shared Google Drive folder for your group this block and upload your .docx file to that folder. Once uploaded, you will all be able to work on that document at one time.\"\n\nReplace-Text `\n    \"In this activity you will all be responsible for the searches and reviewing of items founds. Generally you will all have assigned roles but for in this case you should all be working on the shared document adding different items as explained below. I would like someone to additionally be chosen as a facilitator.\" `\n    \"In this activity, you will all be responsible for conducting searches and reviewing items found. Generally, you will all have assigned roles, but in this case, you should all be working on the shared document, adding different items as explained below. Please choose a facilitator.\"\n\nReplace-Text `\n    \"for this question. In your group I would like you to discuss among yourselves how you would define inequity. Write a one sentence definition.\" `\n    \"for this question. In your group, I would like you to discuss among yourselves how you would define inequity. Write a one-sentence definition.\"\n\nReplace-Text `\n    \"Now you can start to use the internet. Each group member should spend a few minutes looking for different definitions of inequity and examples. Each group member should select a different example. Considering looking for something you have a personal connection to or have thought about before. Include your examples below (make sure to cite your sources, websites are fine in this case).\" `\n    \"Now you can start to use the internet. Each group member should spend a few minutes looking for different definitions of inequity and examples. Each group member should select a different example. Consider looking for something you have a personal connection to or have thought about before. Include your examples below (make sure to cite your sources; websites are fine in this case).\"\n\nReplace-Text `\n    \"Now that you have a definition for inequity you can start too look for it in Iowa. In groups I want you to look for some anecdotal evidence. This type of evidence is primarily what Waypoint has founded their suspicions of inequity on. We are specifically looking for inequity in Iowa and how it relates to having a home.\" `\n    \"Now that you have a definition for inequity, you can start to look for it in Iowa. In groups, please look for some anecdotal evidence. This type of evidence is primarily what Waypoint has founded their suspicions of inequity on. We are specifically looking for inequity in Iowa and how it relates to having a home.\"\n"}
